$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.761.47"
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("D3").Value = "3.008.54"
$ws.Range("E3").Value = "  +3.12%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.19"
$ws.Range("E5").Value = "  +7.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.63"
$ws.Range("E6").Value = "  +8.72%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.434"
$ws.Range("E8").Value = "  +5.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.57"
$ws.Range("E9").Value = "  +11.74%  "
$ws.Range("E10").Value = "  +10.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.356"
$ws.Range("E11").Value = "  +4.59%  "
$ws.Range("E12").Value = "  +3.81%  "
$ws.Range("D13").Value = "3.520.51"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.72"
$ws.Range("E14").Value = "  +8.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000156"
$ws.Range("E15").Value = "  +15.14%  "
$ws.Range("D16").Value = "56.791.33"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").Value = "3.001.38"
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("E18").Value = "  +8.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  +6.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.87"
$ws.Range("E20").Value = "  +8.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.76"
$ws.Range("E21").Value = "  +7.97%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.484"
$ws.Range("E23").Value = "  +7.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.02"
$ws.Range("E24").Value = "  +6.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.175"
$ws.Range("E25").Value = "  +13.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "0.0₃0916"
$ws.Range("E27").Value = "  +11.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.75"
$ws.Range("E28").Value = "  +6.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.09"
$ws.Range("E29").Value = "  +11.40%  "
$ws.Range("E30").Value = "  +11.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.80"
$ws.Range("E31").Value = "  +8.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.70"
$ws.Range("E32").Value = "  +9.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.25"
$ws.Range("E33").Value = "  +6.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("E34").Value = "  +8.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.69"
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("E37").Value = "  +8.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.36"
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("D39").Value = "3.040.17"
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.98"
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +4.74%  "
$ws.Range("D43").Value = "2.272.98"
$ws.Range("E43").Value = "  +9.76%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.42"
$ws.Range("E44").Value = "  +5.40%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +2.92%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.67"
$ws.Range("E46").Value = "  +6.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("E47").Value = "  +22.57%  "
$ws.Range("E48").Value = "  +8.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.83"
$ws.Range("E49").Value = "  +6.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.49"
$ws.Range("E50").Value = "  +7.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0872"
$ws.Range("E51").Value = "  +8.50%  "
